$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J28").Value = 8896.333000000001
$ws.Range("M28").Value = -1339.25
$ws.Range("H28").Value = 3753
$ws.Range("L28").Value = 8896.333000000001
$ws.Range("I28").Value = 1824.25
$ws.Range("K28").Value = 1824.25
$ws.Range("N28").Value = -9866.333000000001
$ws.Range("I40").Value = 5507.4287
$ws.Range("K40").Value = 5507.4287
$ws.Range("H40").Value = 7549.857
$ws.Range("M40").Value = -5332.4287
$ws.Range("M41").Value = -9.5
$ws.Range("H41").Value = 1438.2
$ws.Range("N41").Value = -2677.7273
$ws.Range("J41").Value = 1797.7273
$ws.Range("L41").Value = 1797.7273
$ws.Range("K41").Value = 449.5
$ws.Range("I41").Value = 449.5
$ws.Range("M74").Value = -8559.666999999999
$ws.Range("K74").Value = 9495.666999999999
$ws.Range("I74").Value = 9495.666999999999
$ws.Range("H74").Value = 9495.666999999999
$ws.Range("K77").Value = 47478.335
$ws.Range("H77").Value = 9495.666999999999
$ws.Range("M77").Value = -42798.335
$ws.Range("I77").Value = 9495.666999999999
$ws.Range("J92").Value = 1521.6666
$ws.Range("I92").Value = 1610.1333
$ws.Range("K92").Value = 1610.1333
$ws.Range("M92").Value = -362.1333
$ws.Range("H92").Value = 1595.3889
$ws.Range("N92").Value = -4017.6666
$ws.Range("L92").Value = 1521.6666
$ws.Range("I113").Value = 4000
$ws.Range("H113").Value = 6189.4
$ws.Range("M113").Value = -746
$ws.Range("K113").Value = 4000
$ws.Range("L117").Value = 103266.664
$ws.Range("J117").Value = 103266.664
$ws.Range("H117").Value = 103266.664
$ws.Range("N117").Value = -112444.664
$ws.Range("J137").Value = 3406.9285
$ws.Range("K137").Value = 6928.5
$ws.Range("I137").Value = 2309.5
$ws.Range("M137").Value = -4378.5
$ws.Range("H137").Value = 3077.7
$ws.Range("L137").Value = 10220.7855
$ws.Range("N137").Value = -15320.7855

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K45").Value = 3445.5
$ws.Range("M45").Value = -3068.5
$ws.Range("H45").Value = 3828.7334
$ws.Range("I45").Value = 3445.5
$ws.Range("H63").Value = 13632.667
$ws.Range("N63").Value = -17571
$ws.Range("M63").Value = -115
$ws.Range("K63").Value = 801
$ws.Range("I63").Value = 801
$ws.Range("J63").Value = 16199
$ws.Range("L63").Value = 16199
$ws.Range("K66").Value = 4005
$ws.Range("N66").Value = -87859
$ws.Range("M66").Value = -573
$ws.Range("J66").Value = 16199
$ws.Range("L66").Value = 80995
$ws.Range("H66").Value = 13632.667
$ws.Range("I66").Value = 801
$ws.Range("I122").Value = 4467.25
$ws.Range("K122").Value = 13401.75
$ws.Range("H122").Value = 4573.2
$ws.Range("M122").Value = -10951.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2153.9
$ws.Range("M22").Value = -2019.625
$ws.Range("I22").Value = 2192.625
$ws.Range("K22").Value = 2192.625

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K16").Value = 974.75
$ws.Range("H16").Value = 974.75
$ws.Range("I16").Value = 974.75
$ws.Range("M16").Value = -687.75
$ws.Range("H22").Value = 238
$ws.Range("M22").Value = 127.5
$ws.Range("I22").Value = 222.5
$ws.Range("K22").Value = 222.5
$ws.Range("K58").Value = 1695.4546
$ws.Range("H58").Value = 2367.4285
$ws.Range("I58").Value = 1695.4546
$ws.Range("M58").Value = -1492.4546
$ws.Range("H105").Value = 1693.2858
$ws.Range("I105").Value = 1145.5714
$ws.Range("K105").Value = 1145.5714
$ws.Range("M105").Value = 601.4286
$ws.Range("I113").Value = 974.75
$ws.Range("H113").Value = 974.75
$ws.Range("M113").Value = 1195.25
$ws.Range("K113").Value = 974.75
$ws.Range("I122").Value = 1887.25
$ws.Range("K122").Value = 5661.75
$ws.Range("H122").Value = 1533.75
$ws.Range("M122").Value = -3211.75
$ws.Range("M136").Value = -2536.3638
$ws.Range("K136").Value = 5086.3638
$ws.Range("H136").Value = 2367.4285
$ws.Range("I136").Value = 1695.4546

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L4").Value = 1503.6
$ws.Range("N4").Value = -1727.6
$ws.Range("J4").Value = 501.2
$ws.Range("H4").Value = 561683.1
$ws.Range("K87").Value = 9021
$ws.Range("M87").Value = -7773
$ws.Range("I87").Value = 3007
$ws.Range("H87").Value = 4338
$ws.Range("K90").Value = 27063
$ws.Range("H90").Value = 4338
$ws.Range("I90").Value = 3007
$ws.Range("M90").Value = -20823
$ws.Range("I92").Value = 449
$ws.Range("K92").Value = 1347
$ws.Range("M92").Value = -99
$ws.Range("H92").Value = 290.5
$ws.Range("K111").Value = 450.75
$ws.Range("M111").Value = 2616.25
$ws.Range("H111").Value = 150.25
$ws.Range("I111").Value = 150.25
$ws.Range("I141").Value = 1169.6666
$ws.Range("L141").Value = 21000
$ws.Range("N141").Value = -31360
$ws.Range("M141").Value = 1671.0002
$ws.Range("J141").Value = 7000
$ws.Range("H141").Value = 3501.8
$ws.Range("K141").Value = 3508.9998

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I2").Value = 31.09091
$ws.Range("J2").Value = 199.83333
$ws.Range("N2").Value = -425.83333
$ws.Range("M2").Value = 81.90908999999999
$ws.Range("K2").Value = 31.09091
$ws.Range("H2").Value = 90.64706
$ws.Range("L2").Value = 199.83333
$ws.Range("L39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("H39").Value = 0
$ws.Range("M43").Value = 123
$ws.Range("I43").Value = 28
$ws.Range("H43").Value = 14949.583
$ws.Range("K43").Value = 28
$ws.Range("K57").Value = 25000
$ws.Range("M57").Value = -24180
$ws.Range("H57").Value = 49191.25
$ws.Range("J57").Value = 73382.5
$ws.Range("L57").Value = 73382.5
$ws.Range("I57").Value = 25000
$ws.Range("N57").Value = -75022.5
$ws.Range("N97").Value = -1425
$ws.Range("J97").Value = 433
$ws.Range("M97").Value = 372
$ws.Range("L97").Value = 433
$ws.Range("I97").Value = 124
$ws.Range("H97").Value = 309.4
$ws.Range("K97").Value = 124
$ws.Range("K102").Value = 4999
$ws.Range("M102").Value = -3377
$ws.Range("H102").Value = 4999
$ws.Range("I102").Value = 4999
$ws.Range("N107").Value = -7651.625
$ws.Range("J107").Value = 3811.625
$ws.Range("L107").Value = 3811.625
$ws.Range("H107").Value = 2234
$ws.Range("L132").Value = 0
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("N132").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 962.5
$ws.Range("N22").Value = -1590
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("H27").Value = 962.5
$ws.Range("L27").Value = 1000
$ws.Range("N27").Value = -1214
$ws.Range("J27").Value = 1000
$ws.Range("J46").Value = 4333.3335
$ws.Range("L46").Value = 4333.3335
$ws.Range("H46").Value = 4333.3335
$ws.Range("N46").Value = -4709.3335
$ws.Range("J68").Value = 10000
$ws.Range("N68").Value = -11498
$ws.Range("H68").Value = 10000
$ws.Range("L68").Value = 10000
$ws.Range("L71").Value = 50000
$ws.Range("N71").Value = -57488
$ws.Range("H71").Value = 10000
$ws.Range("J71").Value = 10000
$ws.Range("H82").Value = 7400
$ws.Range("I82").Value = 0
$ws.Range("N82").Value = -8122
$ws.Range("L82").Value = 7400
$ws.Range("K82").Value = 0
$ws.Range("J82").Value = 7400
$ws.Range("J85").Value = 7400
$ws.Range("N85").Value = -9896
$ws.Range("K85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("H85").Value = 7400
$ws.Range("L85").Value = 7400
$ws.Range("H122").Value = 5515
$ws.Range("I132").Value = 9746
$ws.Range("H132").Value = 15745.333
$ws.Range("K132").Value = 29238
$ws.Range("M132").Value = -26708
$ws.Range("M136").Value = -5276.25
$ws.Range("J136").Value = 4058.6
$ws.Range("K136").Value = 7826.25
$ws.Range("H136").Value = 3267.7727
$ws.Range("L136").Value = 12175.8
$ws.Range("N136").Value = -17275.8
$ws.Range("I136").Value = 2608.75
$ws.Range("M82").ClearContents()
$ws.Range("M85").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J62").Value = 12571.429
$ws.Range("L62").Value = 12571.429
$ws.Range("K62").Value = 8895
$ws.Range("N62").Value = -13819.429
$ws.Range("M62").Value = -8271
$ws.Range("I62").Value = 8895
$ws.Range("H62").Value = 11468.5
$ws.Range("J65").Value = 12571.429
$ws.Range("L65").Value = 62857.145
$ws.Range("H65").Value = 11468.5
$ws.Range("M65").Value = -41355
$ws.Range("N65").Value = -69097.145
$ws.Range("K65").Value = 44475
$ws.Range("I65").Value = 8895
$ws.Range("J100").Value = 549
$ws.Range("N100").Value = -2180
$ws.Range("H100").Value = 1274.5
$ws.Range("L100").Value = 1098
$ws.Range("N107").Value = -7723.5
$ws.Range("J107").Value = 1294.5
$ws.Range("L107").Value = 3883.5
$ws.Range("H107").Value = 3036.4
$ws.Range("M136").Value = -1078.05
$ws.Range("J136").Value = 4558.3
$ws.Range("K136").Value = 3628.05
$ws.Range("H136").Value = 2325.6667
$ws.Range("L136").Value = 13674.9
$ws.Range("N136").Value = -18774.9
$ws.Range("I136").Value = 1209.35

